$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Sheet1" to "Shivam Dube"
$ws.Name = "Shivam Dube"

# Row 1: headers. "matchNo" is a new first column, so every other header shifts one column right.
$headers = @("matchNo", "teamName", "batterName", "states", "runs", "balls", "fours", "sixes", "sr", "opponentTeamName", "venue", "date", "result")
for ($col = 1; $col -le $headers.Length; $col++) {
    $ws.Cells.Item(1, $col).Value = $headers[$col - 1]
}

# Rows 2-10: full match-by-match batting log for Shivam Dube (row 2 is new; the former
# sole data row - "c Dhawan b Avesh Khan" - now lands on row 9 once matchNo is prepended).
$data = @(
    @("47th", "Rajasthan Royals", "Shivam Dube", "", "64", "42", "4", "4", "152.38", "Chennai Super Kings", "Abu Dhabi", "October 02", "Royals won by 7 wickets (with 15 balls remaining)"),
    @("24th", "Rajasthan Royals", "Shivam Dube", "c & b Bumrah", "35", "31", "2", "2", "112.90", "Mumbai Indians", "Delhi", "April 29", "Mumbai won by 7 wickets (with 9 balls remaining)"),
    @("51st", "Rajasthan Royals", "Shivam Dube", "b Neesham", "3", "8", "0", "0", "37.50", "Mumbai Indians", "Sharjah", "October 05", "Mumbai won by 8 wickets (with 70 balls remaining)"),
    @("18th", "Rajasthan Royals", "Shivam Dube", "c Prasidh Krishna b Varun", "22", "18", "2", "1", "122.22", "Kolkata Knight Riders", "Wankhede", "April 24", "Royals won by 6 wickets (with 7 balls remaining)"),
    @("16th", "Rajasthan Royals", "Shivam Dube", "c Maxwell b Richardson", "46", "32", "5", "2", "143.75", "Royal Challengers Bangalore", "Wankhede", "April 22", "RCB won by 10 wickets (with 21 balls remaining)"),
    @("54th", "Rajasthan Royals", "Shivam Dube", "b Shivam Mavi", "18", "20", "0", "1", "90.00", "Kolkata Knight Riders", "Sharjah", "October 07", "KKR won by 86 runs"),
    @("4th", "Rajasthan Royals", "Shivam Dube", "c Hooda b Arshdeep Singh", "23", "15", "3", "0", "153.33", "Punjab Kings", "Wankhede", "April 12", "Punjab Kings won by 4 runs"),
    @("7th", "Rajasthan Royals", "Shivam Dube", "c Dhawan b Avesh Khan", "2", "7", "0", "0", "28.57", "Delhi Capitals", "Wankhede", "April 15", "Royals won by 3 wickets (with 2 balls remaining)"),
    @("12th", "Rajasthan Royals", "Shivam Dube", "lbw b Jadeja", "17", "20", "2", "0", "85.00", "Chennai Super Kings", "Wankhede", "April 19", "Super Kings won by 45 runs")
)

# Columns that hold digit-only text in the source data (states, runs, balls, fours, sixes, sr)
# must stay text, not be reinterpreted as numbers - same as the original sheet, which
# flags the whole range with numberStoredAsText in ignoredErrors.
$textColumns = @(4, 5, 6, 7, 8, 9)

for ($r = 0; $r -lt $data.Length; $r++) {
    $rowValues = $data[$r]
    $rowIndex = $r + 2
    for ($col = 1; $col -le $rowValues.Length; $col++) {
        $cellValue = $rowValues[$col - 1]
        $cell = $ws.Cells.Item($rowIndex, $col)
        if ($textColumns -contains $col) {
            # Leading apostrophe forces a text (quote-prefixed) cell, exactly like typing
            # '64 into Excel, so "64" is kept as text instead of becoming the number 64.
            # Re-normalize the style afterwards so only the cells stored type/value
            # changes, without leaving a quote-prefix number format behind.
            $cell.Value = "'" + $cellValue
            $cell.Style = "Normal"
        } else {
            $cell.Value = $cellValue
        }
    }
}

Write-Output "done"
